$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.223.37"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "'1.662.35"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'217.85"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'0.5227"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.2643"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "'0.06282"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").Value = "'20.79"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "'1.666.48"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'4.429"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "'1.889.42"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "'0.5431"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "'0.0₅8141"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").Value = "'64.50"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "'26.252.59"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'4.644"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "'192.96"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "'10.06"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "'6.038"
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "'139.92"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'0.1229"
$ws.Range("E26").Value = "  -4.62%  "
$ws.Range("D27").Value = "'7.158"
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("D28").Value = "'16.09"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'1.413"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").Value = "'0.06080"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'1.278"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'3.575"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "'3.257"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").Value = "'1.620"
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("D35").Value = "'0.9643"
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").Value = "'2.427"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'2.783"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'0.5672"
$ws.Range("E38").Value = "  -8.28%  "
$ws.Range("D39").Value = "'0.01591"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").Value = "'5.977"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").Value = "'0.8560"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").Value = "'1.012.52"
$ws.Range("E43").Value = "  -7.38%  "
$ws.Range("D44").Value = "'100.23"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'1.804.89"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.06"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.0₈107"
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("D48").Value = "'1.010"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "'8.004"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05181"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.475"
$ws.Range("E51").Value = "  -1.20%  "
